# ex07 - p3 - matvec runs single threaded, some corrections in p2
#
# This script edits the "P2" worksheet of the workbook:
#  - corrects the measured GPU (col B) and CPU (col D) timings
#  - adds a new "rel CPU/GPU" column (E) with a =D/B formula
#  - adds a small helper table (F21:G30) with vector size / CPU time
#  - repositions the chart and its legend to match the new data extents
#  - updates the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("P2")

# ---------------------------------------------------------------------
# 1. Corrected measurement values for "OpenCL on GPU" (col B) and
#    "OpenCL on CPU" (col D), rows 5..13 (vector sizes 1000..10000000)
# ---------------------------------------------------------------------
$ws.Range("B5").Value  = 0.000038
$ws.Range("B6").Value  = 0.000039
$ws.Range("B7").Value  = 0.00004
$ws.Range("B8").Value  = 0.000046
$ws.Range("B9").Value  = 0.000039
$ws.Range("B10").Value = 0.000052
$ws.Range("B11").Value = 0.000106
$ws.Range("B12").Value = 0.000177
$ws.Range("B13").Value = 0.000499

$ws.Range("D5").Value  = 0.000393
$ws.Range("D6").Value  = 0.000394
$ws.Range("D7").Value  = 0.000415
$ws.Range("D8").Value  = 0.000485
$ws.Range("D9").Value  = 0.000449
$ws.Range("D10").Value = 0.000497
$ws.Range("D11").Value = 0.001994
$ws.Range("D12").Value = 0.006804
$ws.Range("D13").Value = 0.072691

# ---------------------------------------------------------------------
# 2. New "rel CPU/GPU" column
# ---------------------------------------------------------------------
$ws.Range("E4").Value = "rel CPU/GPU"

$ws.Range("E5").Formula = "=D5/B5"
$ws.Range("E6:E13").Formula = "=D6/B6"

# ---------------------------------------------------------------------
# 3. Helper table F21:G30 (vector size / CPU time), formatted like the
#    existing A/D and B/C helper columns
# ---------------------------------------------------------------------
$ws.Range("A21").Copy()
$ws.Range("F21").PasteSpecial(-4122)
$ws.Range("B5").Copy()
$ws.Range("G21").PasteSpecial(-4122)

$sizes = @(1000, 3000, 10000, 30000, 100000, 300000, 1000000, 3000000, 10000000)
$cpu   = @(0.000393, 0.000394, 0.000415, 0.000485, 0.000449, 0.000497, 0.001994, 0.006804, 0.072691)

for ($i = 0; $i -lt $sizes.Length; $i++) {
    $r = 22 + $i
    $ws.Range("A22").Copy()
    $ws.Range("F$r").PasteSpecial(-4122)
    $ws.Range("B5").Copy()
    $ws.Range("G$r").PasteSpecial(-4122)
    $ws.Range("F$r").Value = $sizes[$i]
    $ws.Range("G$r").Value = $cpu[$i]
}

# ---------------------------------------------------------------------
# 4. Reposition chart & legend to match new layout
# ---------------------------------------------------------------------
$chartObj = $ws.ChartObjects().Item(1)
$chartObj.Left = 317.5323031496063
$chartObj.Top = 36.76291338582677
$chartObj.Width = 443.49999999999994
$chartObj.Height = 216.0

$legend = $chartObj.Chart.Legend
$legend.Left = 0.17977680354216358
$legend.Top = 0.056133712452610104
$legend.Width = 0.32968416871976269
$legend.Height = 0.27451698745990083

# ---------------------------------------------------------------------
# 5. Update active selection
# ---------------------------------------------------------------------
$ws.Range("E17").Select()
